{"js": "// Update the three-digit-by-one-digit multiplication answers in the table.\n// Each original equation string is unique in the document, so a targeted\n// search + replace for each pair is unambiguous and safe.\n\nconst replacements = [\n  [\"238\u00d79=2142\", \"850\u00d78=6800\"],\n  [\"360\u00d75=1800\", \"145\u00d76=870\"],\n  [\"815\u00d72=1630\", \"737\u00d74=2948\"],\n  [\"425\u00d72=850\", \"949\u00d72=1898\"],\n  [\"522\u00d77=3654\", \"964\u00d73=2892\"],\n  [\"851\u00d72=1702\", \"927\u00d76=5562\"],\n  [\"257\u00d78=2056\", \"358\u00d77=2506\"],\n  [\"493\u00d74=1972\", \"537\u00d72=1074\"],\n  [\"870\u00d73=2610\", \"154\u00d74=616\"],\n  [\"493\u00d76=2958\", \"281\u00d77=1967\"],\n  [\"945\u00d72=1890\", \"975\u00d72=1950\"],\n  [\"750\u00d76=4500\", \"997\u00d79=8973\"],\n  [\"658\u00d77=4606\", \"363\u00d76=2178\"],\n  [\"972\u00d75=4860\", \"631\u00d73=1893\"],\n  [\"700\u00d73=2100\", \"578\u00d79=5202\"],\n  [\"791\u00d77=5537\", \"626\u00d74=2504\"],\n  [\"537\u00d74=2148\", \"920\u00d72=1840\"],\n  [\"412\u00d76=2472\", \"152\u00d73=456\"],\n  [\"732\u00d78=5856\", \"976\u00d78=7808\"],\n  [\"652\u00d78=5216\", \"864\u00d74=3456\"],\n  [\"745\u00d72=1490\", \"285\u00d79=2565\"],\n  [\"765\u00d75=3825\", \"376\u00d75=1880\"],\n  [\"264\u00d76=1584\", \"270\u00d76=1620\"],\n  [\"548\u00d75=2740\", \"296\u00d74=1184\"],\n  [\"606\u00d72=1212\", \"878\u00d78=7024\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit-by-one-digit multiplication answers in the table.\n# Each original equation string is unique in the document, so a targeted\n# Find/Replace (wdReplaceOne) for each pair is unambiguous and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"238\u00d79=2142\", \"850\u00d78=6800\"),\n    @(\"360\u00d75=1800\", \"145\u00d76=870\"),\n    @(\"815\u00d72=1630\", \"737\u00d74=2948\"),\n    @(\"425\u00d72=850\",  \"949\u00d72=1898\"),\n    @(\"522\u00d77=3654\", \"964\u00d73=2892\"),\n    @(\"851\u00d72=1702\", \"927\u00d76=5562\"),\n    @(\"257\u00d78=2056\", \"358\u00d77=2506\"),\n    @(\"493\u00d74=1972\", \"537\u00d72=1074\"),\n    @(\"870\u00d73=2610\", \"154\u00d74=616\"),\n    @(\"493\u00d76=2958\", \"281\u00d77=1967\"),\n    @(\"945\u00d72=1890\", \"975\u00d72=1950\"),\n    @(\"750\u00d76=4500\", \"997\u00d79=8973\"),\n    @(\"658\u00d77=4606\", \"363\u00d76=2178\"),\n    @(\"972\u00d75=4860\", \"631\u00d73=1893\"),\n    @(\"700\u00d73=2100\", \"578\u00d79=5202\"),\n    @(\"791\u00d77=5537\", \"626\u00d74=2504\"),\n    @(\"537\u00d74=2148\", \"920\u00d72=1840\"),\n    @(\"412\u00d76=2472\", \"152\u00d73=456\"),\n    @(\"732\u00d78=5856\", \"976\u00d78=7808\"),\n    @(\"652\u00d78=5216\", \"864\u00d74=3456\"),\n    @(\"745\u00d72=1490\", \"285\u00d79=2565\"),\n    @(\"765\u00d75=3825\", \"376\u00d75=1880\"),\n    @(\"264\u00d76=1584\", \"270\u00d76=1620\"),\n    @(\"548\u00d75=2740\", \"296\u00d74=1184\"),\n    @(\"606\u00d72=1212\", \"878\u00d78=7024\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$new, 2)\n}\n"}
